$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "23.353.73"
Set-TextValue "E2" "  -0.34%  "
Set-TextValue "D3" "1.627.21"
Set-TextValue "E3" "  -0.84%  "
Set-TextValue "E4" "  +0.11%  "
Set-TextValue "D5" "1.000"
Set-TextValue "E5" "  +0.00%  "
Set-TextValue "D6" "301.86"
Set-TextValue "E6" "  -1.05%  "
Set-TextValue "D7" "0.3754"
Set-TextValue "E7" "  +0.49%  "
Set-TextValue "D8" "0.3626"
Set-TextValue "E8" "  -0.15%  "
Set-TextValue "D9" "51.49"
Set-TextValue "E9" "  -1.58%  "
Set-TextValue "D10" "0.08166"
Set-TextValue "E10" "  +0.52%  "
Set-TextValue "D11" "1.217"
Set-TextValue "E11" "  -2.88%  "
Set-TextValue "E12" "  +0.10%  "
Set-TextValue "E13" "  -2.88%  "
Set-TextValue "E14" "  -2.18%  "
Set-TextValue "E15" "  -2.72%  "
Set-TextValue "D16" "7.286"
Set-TextValue "E16" "  -0.02%  "
Set-TextValue "D17" "1.598.45"
Set-TextValue "E17" "  -2.06%  "
Set-TextValue "D18" "94.60"
Set-TextValue "E18" "  +0.18%  "
Set-TextValue "D19" "0.06948"
Set-TextValue "E19" "  +0.71%  "
Set-TextValue "D20" "17.53"
Set-TextValue "E20" "  -3.35%  "
Set-TextValue "D21" "6.552"
Set-TextValue "E21" "  +0.64%  "
Set-TextValue "D22" "1.000"
Set-TextValue "E22" "  -0.03%  "
Set-TextValue "E23" "  -2.27%  "
Set-TextValue "D24" "23.356.77"
Set-TextValue "E24" "  -0.38%  "
Set-TextValue "D25" "2.489"
Set-TextValue "E25" "  +3.18%  "
Set-TextValue "D26" "3.066"
Set-TextValue "E26" "  -0.69%  "
Set-TextValue "D27" "21.12"
Set-TextValue "E27" "  -0.38%  "
Set-TextValue "D28" "150.36"
Set-TextValue "E28" "  -0.84%  "
Set-TextValue "D29" "5.264"
Set-TextValue "E29" "  -1.34%  "
Set-TextValue "D30" "132.75"
Set-TextValue "E30" "  -2.14%  "
Set-TextValue "D31" "1.797.87"
Set-TextValue "E31" "  -0.74%  "
Set-TextValue "D32" "6.598"
Set-TextValue "E32" "  -3.37%  "
Set-TextValue "D33" "2.156"
Set-TextValue "E33" "  -5.39%  "
Set-TextValue "D34" "1.060"
Set-TextValue "E34" "  +11.42%  "
Set-TextValue "D35" "11.18"
Set-TextValue "E35" "  +7.89%  "
Set-TextValue "D36" "0.02755"
Set-TextValue "E36" "  -1.97%  "
Set-TextValue "D37" "0.2486"
Set-TextValue "E37" "  -1.44%  "
Set-TextValue "D38" "0.08746"
Set-TextValue "E38" "  -0.32%  "
Set-TextValue "D39" "0.07121"
Set-TextValue "E39" "  -1.57%  "
Set-TextValue "E40" "  -2.51%  "
Set-TextValue "D41" "0.6963"
Set-TextValue "E41" "  -1.44%  "
Set-TextValue "D42" "1.325"
Set-TextValue "E42" "  -3.35%  "
Set-TextValue "D43" "15.74"
Set-TextValue "E43" "  -2.16%  "
Set-TextValue "D44" "11.95"
Set-TextValue "E44" "  -4.33%  "
Set-TextValue "D45" "0.6431"
Set-TextValue "E45" "  -1.47%  "
Set-TextValue "D46" "0.9996"
Set-TextValue "E46" "  +0.04%  "
Set-TextValue "D47" "2.268"
Set-TextValue "E47" "  -2.64%  "
Set-TextValue "D48" "3.955"
Set-TextValue "E48" "  -1.38%  "
Set-TextValue "D49" "0.07964"
Set-TextValue "E49" "  -0.08%  "
Set-TextValue "D50" "126.96"
Set-TextValue "E50" "  -1.25%  "
Set-TextValue "D51" "1.187"
Set-TextValue "E51" "  -1.06%  "
